$wb = $excel.ActiveWorkbook

$xlToLeft = -4159

# A new student, "alumno6@mail.com" / Quique Quiroga, attended none of the
# sessions so far: add him as a new row (row 6) at the bottom of every
# monthly attendance sheet, marking every date column as absent ("A").
foreach ($ws in $wb.Worksheets) {
    $newRow = 6
    $lastCol = $ws.Cells.Item(1, $ws.Columns.Count).End($xlToLeft).Column

    $ws.Cells.Item($newRow, 1).Value = "alumno6@mail.com"
    $ws.Cells.Item($newRow, 2).Value = "Quique"
    $ws.Cells.Item($newRow, 3).Value = "Quiroga"

    for ($col = 4; $col -le $lastCol; $col++) {
        $ws.Cells.Item($newRow, $col).Value = "A"
    }
}

# In "Julio" everybody actually showed up on 5/7/22 (column E): mark that
# whole column, including the new student, as present ("P").
$julio = $wb.Worksheets.Item("Julio")
$julio.Range("E2:E6").Value = "P"
